# Auto-generated edit script applying the crypto price/volume update
# described by the commit 'Updated cryptos list on Mon Jul  8 19:31:08 UTC 2024 with GitHub Actions'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.358.58'
$ws.Range('E2').Value = '  -1.25%  '

$ws.Range('D3').Value = '2.999.53'
$ws.Range('E3').Value = '  +0.50%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '507.02'
$ws.Range('E5').Value = '  +1.38%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.15'
$ws.Range('E6').Value = '  +0.42%  '

$ws.Range('E7').Value = '  +0.02%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.431'
$ws.Range('E8').Value = '  +0.55%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.13'
$ws.Range('E9').Value = '  -2.05%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.108'
$ws.Range('E10').Value = '  -0.03%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.367'
$ws.Range('E11').Value = '  +2.80%  '

$ws.Range('D12').Value = '3.507.77'
$ws.Range('E12').Value = '  +0.35%  '

$ws.Range('E13').Value = '  -0.67%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.66'
$ws.Range('E14').Value = '  -1.31%  '

$ws.Range('E15').Value = '  +2.26%  '

$ws.Range('D16').Value = '56.354.75'
$ws.Range('E16').Value = '  -1.33%  '

$ws.Range('D17').Value = '2.987.85'
$ws.Range('E17').Value = '  +0.30%  '

$ws.Range('E18').Value = '  -1.25%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.95'
$ws.Range('E19').Value = '  +2.69%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.08'
$ws.Range('E20').Value = '  +2.66%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '331.93'
$ws.Range('E21').Value = '  +3.62%  '

$ws.Range('E22').Value = '  +0.20%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.497'
$ws.Range('E23').Value = '  +1.23%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.98'
$ws.Range('E24').Value = '  +3.13%  '

$ws.Range('D25').Value = '3.116.68'
$ws.Range('E25').Value = '  +0.22%  '

$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').Value = '0.0₃0946'
$ws.Range('E26').Value = '  +6.08%  '

$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.165'
$ws.Range('E27').Value = '  +1.21%  '

$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.996'
$ws.Range('E28').Value = '  -0.31%  '

$ws.Range('E29').Value = '  -3.70%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.93'
$ws.Range('E30').Value = '  -2.67%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.79'
$ws.Range('E31').Value = '  +0.73%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.31'
$ws.Range('E32').Value = '  +0.80%  '

$ws.Range('E33').Value = '  -0.45%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '153.74'
$ws.Range('E34').Value = '  -0.50%  '

$ws.Range('E35').Value = '  -1.95%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.85'
$ws.Range('E36').Value = '  +1.19%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '26.28'
$ws.Range('E37').Value = '  +8.08%  '

$ws.Range('E38').Value = '  +0.66%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0662'
$ws.Range('E39').Value = '  -0.32%  '

$ws.Range('D40').Value = '3.032.67'
$ws.Range('E40').Value = '  +0.59%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.71'
$ws.Range('E41').Value = '  -2.85%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.02%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.80'
$ws.Range('E43').Value = '  +1.60%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.652'
$ws.Range('E44').Value = '  +1.20%  '

$ws.Range('D45').Value = '2.181.98'
$ws.Range('E45').Value = '  -0.42%  '

$ws.Range('E46').Value = '  -2.28%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.88'
$ws.Range('E47').Value = '  -1.39%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.928'
$ws.Range('E48').Value = '  -0.65%  '

$ws.Range('E49').Value = '  +1.06%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.47'
$ws.Range('E50').Value = '  +1.48%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0854'
$ws.Range('E51').Value = '  -1.61%  '
